# Auto-generated script applying row permutation of species-occurrence data
# for rows 47-65 (A,B,D,E,F,G,H,I,J,K,L,M,N,Q,R,S,AC,AF columns move;
# C,P,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY columns stay fixed per row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47  <-  source row 49
$ws.Range("A47").Value = 111736257
$ws.Range("B47").Value = 77515
$ws.Range("D47").Value = 'NT'
$ws.Range("E47").Value = 6425
$ws.Range("F47").Value = 'Garnlav'
$ws.Range("G47").Value = 'Alectoria sarmentosa'
$ws.Range("H47").Value = '(Ach.) Ach.'
$ws.Range("I47").Value = ''
$ws.Range("J47").Value = ''
$ws.Range("K47").Value = ''
$ws.Range("N47").Value = ''
$ws.Range("Q47").Value = 616308.8236423519
$ws.Range("R47").Value = 7268903.133137755
$ws.Range("S47").Value = 10
$ws.Range("AF47").Value = ''

# Row 48  <-  source row 50
$ws.Range("A48").Value = 111736402
$ws.Range("B48").Value = 89423
$ws.Range("D48").Value = 'NT'
$ws.Range("E48").Value = 5432
$ws.Range("F48").Value = 'Granticka'
$ws.Range("G48").Value = 'Porodaedalea chrysoloma'
$ws.Range("H48").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("I48").Value = ''
$ws.Range("J48").Value = ''
$ws.Range("K48").Value = ''
$ws.Range("L48").ClearContents()
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = ''
$ws.Range("Q48").Value = 616333.1117616051
$ws.Range("R48").Value = 7268857.179896916
$ws.Range("S48").Value = 10
$ws.Range("AF48").Value = ''

# Row 49  <-  source row 52
$ws.Range("A49").Value = 111736506
$ws.Range("B49").Value = 56398
$ws.Range("D49").Value = 'NT'
$ws.Range("E49").Value = 100109
$ws.Range("F49").Value = 'Tretåig hackspett'
$ws.Range("G49").Value = 'Picoides tridactylus'
$ws.Range("H49").Value = '(Linnaeus, 1758)'
$ws.Range("I49").Value = ''
$ws.Range("J49").ClearContents()
$ws.Range("K49").Value = ''
$ws.Range("L49").Value = ''
$ws.Range("M49").Value = 'födosökande'
$ws.Range("N49").Value = ''
$ws.Range("Q49").Value = 616358.6131022752
$ws.Range("R49").Value = 7268822.486957001
$ws.Range("S49").Value = 25
$ws.Range("AC49").Value = 'Födosökande i gransumpskog. Observerades i över en timme'
$ws.Range("AF49").ClearContents()

# Row 50  <-  source row 48
$ws.Range("A50").Value = 111736370
$ws.Range("B50").Value = 56398
$ws.Range("D50").Value = 'NT'
$ws.Range("E50").Value = 100109
$ws.Range("F50").Value = 'Tretåig hackspett'
$ws.Range("G50").Value = 'Picoides tridactylus'
$ws.Range("H50").Value = '(Linnaeus, 1758)'
$ws.Range("I50").Value = ''
$ws.Range("J50").ClearContents()
$ws.Range("K50").Value = ''
$ws.Range("L50").Value = ''
$ws.Range("M50").Value = 'färska spår'
$ws.Range("N50").Value = ''
$ws.Range("Q50").Value = 616327.1020967637
$ws.Range("R50").Value = 7268872.304318298
$ws.Range("S50").Value = 10
$ws.Range("AF50").ClearContents()

# Row 51  <-  source row 47
$ws.Range("A51").Value = 111736268
$ws.Range("B51").Value = 89401
$ws.Range("D51").Value = 'NT'
$ws.Range("E51").Value = 1108
$ws.Range("F51").Value = 'Harticka'
$ws.Range("G51").Value = 'Pelloporus leporinus'
$ws.Range("H51").Value = '(Fr.) Krieglst.'
$ws.Range("I51").Value = ''
$ws.Range("J51").Value = ''
$ws.Range("K51").Value = ''
$ws.Range("N51").Value = ''
$ws.Range("Q51").Value = 616308.8236423519
$ws.Range("R51").Value = 7268903.133137755
$ws.Range("S51").Value = 10
$ws.Range("AF51").Value = ''

# Row 52  <-  source row 51
$ws.Range("A52").Value = 111736405
$ws.Range("B52").Value = 77515
$ws.Range("D52").Value = 'NT'
$ws.Range("E52").Value = 6425
$ws.Range("F52").Value = 'Garnlav'
$ws.Range("G52").Value = 'Alectoria sarmentosa'
$ws.Range("H52").Value = '(Ach.) Ach.'
$ws.Range("I52").Value = ''
$ws.Range("J52").Value = ''
$ws.Range("K52").Value = ''
$ws.Range("L52").ClearContents()
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = ''
$ws.Range("Q52").Value = 616333.1117616051
$ws.Range("R52").Value = 7268857.179896916
$ws.Range("S52").Value = 10
$ws.Range("AC52").ClearContents()
$ws.Range("AF52").Value = ''

# Row 53  <-  source row 55
$ws.Range("A53").Value = 111778163
$ws.Range("B53").Value = 56398
$ws.Range("D53").Value = 'NT'
$ws.Range("E53").Value = 100109
$ws.Range("F53").Value = 'Tretåig hackspett'
$ws.Range("G53").Value = 'Picoides tridactylus'
$ws.Range("H53").Value = '(Linnaeus, 1758)'
$ws.Range("I53").Value = ''
$ws.Range("K53").Value = ''
$ws.Range("L53").Value = ''
$ws.Range("M53").Value = 'färska spår'
$ws.Range("N53").Value = ''
$ws.Range("Q53").Value = 616207.2556492372
$ws.Range("R53").Value = 7268635.7870906
$ws.Range("S53").Value = 10
$ws.Range("AC53").Value = 'Skalad gran'

# Row 54  <-  source row 63
$ws.Range("A54").Value = 111777411
$ws.Range("B54").Value = 56398
$ws.Range("D54").Value = 'NT'
$ws.Range("E54").Value = 100109
$ws.Range("F54").Value = 'Tretåig hackspett'
$ws.Range("G54").Value = 'Picoides tridactylus'
$ws.Range("H54").Value = '(Linnaeus, 1758)'
$ws.Range("I54").Value = ''
$ws.Range("J54").ClearContents()
$ws.Range("K54").Value = ''
$ws.Range("L54").Value = ''
$ws.Range("M54").Value = 'färska spår'
$ws.Range("N54").Value = ''
$ws.Range("Q54").Value = 616367.7277224116
$ws.Range("R54").Value = 7268802.503264537
$ws.Range("S54").Value = 10
$ws.Range("AC54").Value = 'Skalade granstammar'
$ws.Range("AF54").ClearContents()

# Row 55  <-  source row 64
$ws.Range("A55").Value = 111777499
$ws.Range("B55").Value = 78107
$ws.Range("D55").Value = 'NT'
$ws.Range("E55").Value = 6453
$ws.Range("F55").Value = 'Vedskivlav'
$ws.Range("G55").Value = 'Hertelidea botryosa'
$ws.Range("H55").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("I55").Value = ''
$ws.Range("J55").Value = ''
$ws.Range("K55").Value = ''
$ws.Range("L55").ClearContents()
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = ''
$ws.Range("Q55").Value = 616426.5202303537
$ws.Range("R55").Value = 7268746.301918368
$ws.Range("S55").Value = 10
$ws.Range("AC55").ClearContents()
$ws.Range("AF55").Value = ''

# Row 56  <-  source row 62
$ws.Range("A56").Value = 111778126
$ws.Range("B56").Value = 89405
$ws.Range("D56").Value = 'NT'
$ws.Range("E56").Value = 1202
$ws.Range("F56").Value = 'Ullticka'
$ws.Range("G56").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H56").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I56").Value = ''
$ws.Range("J56").Value = ''
$ws.Range("K56").Value = ''
$ws.Range("N56").Value = ''
$ws.Range("Q56").Value = 616202.3044715263
$ws.Range("R56").Value = 7268603.611313918
$ws.Range("S56").Value = 10
$ws.Range("AF56").Value = ''

# Row 57  <-  source row 53
$ws.Range("A57").Value = 111778248
$ws.Range("B57").Value = 56398
$ws.Range("D57").Value = 'NT'
$ws.Range("E57").Value = 100109
$ws.Range("F57").Value = 'Tretåig hackspett'
$ws.Range("G57").Value = 'Picoides tridactylus'
$ws.Range("H57").Value = '(Linnaeus, 1758)'
$ws.Range("I57").Value = ''
$ws.Range("K57").Value = ''
$ws.Range("L57").Value = ''
$ws.Range("M57").Value = 'färsk spillning'
$ws.Range("N57").Value = ''
$ws.Range("Q57").Value = 616162.9874832245
$ws.Range("R57").Value = 7268630.281087617
$ws.Range("S57").Value = 10
$ws.Range("AC57").Value = 'Skalad gran'

# Row 58  <-  source row 60
$ws.Range("A58").Value = 111777940
$ws.Range("B58").Value = 90678
$ws.Range("D58").Value = 'LC'
$ws.Range("E58").Value = 4366
$ws.Range("F58").Value = 'Skarp dropptaggsvamp'
$ws.Range("G58").Value = 'Hydnellum peckii'
$ws.Range("H58").Value = 'Banker'
$ws.Range("I58").Value = ''
$ws.Range("J58").Value = ''
$ws.Range("K58").Value = ''
$ws.Range("N58").Value = ''
$ws.Range("Q58").Value = 616438.7745429112
$ws.Range("R58").Value = 7268803.685732875
$ws.Range("S58").Value = 25
$ws.Range("AF58").Value = ''

# Row 59  <-  source row 54
$ws.Range("A59").Value = 111777494
$ws.Range("B59").Value = 90854
$ws.Range("D59").Value = 'NT'
$ws.Range("E59").Value = 2079
$ws.Range("F59").Value = 'Nordtagging'
$ws.Range("G59").Value = 'Odonticium romellii'
$ws.Range("H59").Value = '(S.Lundell) Parmasto'
$ws.Range("I59").Value = ''
$ws.Range("J59").Value = ''
$ws.Range("K59").Value = ''
$ws.Range("N59").Value = ''
$ws.Range("Q59").Value = 616426.5202303537
$ws.Range("R59").Value = 7268746.301918368
$ws.Range("S59").Value = 10
$ws.Range("AF59").Value = ''

# Row 60  <-  source row 59
$ws.Range("A60").Value = 111777467
$ws.Range("B60").Value = 89423
$ws.Range("D60").Value = 'NT'
$ws.Range("E60").Value = 5432
$ws.Range("F60").Value = 'Granticka'
$ws.Range("G60").Value = 'Porodaedalea chrysoloma'
$ws.Range("H60").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("I60").Value = ''
$ws.Range("J60").Value = ''
$ws.Range("K60").Value = ''
$ws.Range("N60").Value = ''
$ws.Range("Q60").Value = 616413.4864248879
$ws.Range("R60").Value = 7268760.315060399
$ws.Range("S60").Value = 10
$ws.Range("AF60").Value = ''

# Row 61  <-  source row 58
$ws.Range("A61").Value = 111777447
$ws.Range("B61").Value = 89405
$ws.Range("D61").Value = 'NT'
$ws.Range("E61").Value = 1202
$ws.Range("F61").Value = 'Ullticka'
$ws.Range("G61").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H61").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I61").Value = ''
$ws.Range("J61").Value = ''
$ws.Range("K61").Value = ''
$ws.Range("N61").Value = ''
$ws.Range("Q61").Value = 616379.7321599644
$ws.Range("R61").Value = 7268803.814155157
$ws.Range("S61").Value = 10
$ws.Range("AF61").Value = ''

# Row 62  <-  source row 56
$ws.Range("A62").Value = 111777331
$ws.Range("B62").Value = 89423
$ws.Range("D62").Value = 'NT'
$ws.Range("E62").Value = 5432
$ws.Range("F62").Value = 'Granticka'
$ws.Range("G62").Value = 'Porodaedalea chrysoloma'
$ws.Range("H62").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("I62").Value = ''
$ws.Range("J62").Value = ''
$ws.Range("K62").Value = ''
$ws.Range("N62").Value = ''
$ws.Range("Q62").Value = 616362.7639770868
$ws.Range("R62").Value = 7268822.653031595
$ws.Range("S62").Value = 10
$ws.Range("AF62").Value = ''

# Row 63  <-  source row 65
$ws.Range("A63").Value = 111777491
$ws.Range("B63").Value = 56398
$ws.Range("D63").Value = 'NT'
$ws.Range("E63").Value = 100109
$ws.Range("F63").Value = 'Tretåig hackspett'
$ws.Range("G63").Value = 'Picoides tridactylus'
$ws.Range("H63").Value = '(Linnaeus, 1758)'
$ws.Range("I63").Value = ''
$ws.Range("K63").Value = ''
$ws.Range("L63").Value = ''
$ws.Range("M63").Value = 'färska spår'
$ws.Range("N63").Value = ''
$ws.Range("Q63").Value = 616426.5202303537
$ws.Range("R63").Value = 7268746.301918368
$ws.Range("S63").Value = 10
$ws.Range("AC63").Value = 'Skalade granstammar'

# Row 64  <-  source row 57
$ws.Range("A64").Value = 111777380
$ws.Range("B64").Value = 56398
$ws.Range("D64").Value = 'NT'
$ws.Range("E64").Value = 100109
$ws.Range("F64").Value = 'Tretåig hackspett'
$ws.Range("G64").Value = 'Picoides tridactylus'
$ws.Range("H64").Value = '(Linnaeus, 1758)'
$ws.Range("I64").Value = ''
$ws.Range("J64").ClearContents()
$ws.Range("K64").Value = ''
$ws.Range("L64").Value = ''
$ws.Range("M64").Value = 'färska spår'
$ws.Range("N64").Value = ''
$ws.Range("Q64").Value = 616414.0528149965
$ws.Range("R64").Value = 7268860.418718725
$ws.Range("S64").Value = 10
$ws.Range("AC64").Value = 'Skalade stammar'
$ws.Range("AF64").ClearContents()

# Row 65  <-  source row 61
$ws.Range("A65").Value = 111778005
$ws.Range("B65").Value = 89369
$ws.Range("D65").Value = 'LC'
$ws.Range("E65").Value = 5447
$ws.Range("F65").Value = 'Vedticka'
$ws.Range("G65").Value = 'Fuscoporia viticola'
$ws.Range("H65").Value = '(Schwein.) Murrill'
$ws.Range("I65").Value = ''
$ws.Range("J65").Value = ''
$ws.Range("K65").Value = ''
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = ''
$ws.Range("Q65").Value = 616499.3130462242
$ws.Range("R65").Value = 7268610.508796399
$ws.Range("S65").Value = 10
$ws.Range("AC65").ClearContents()
$ws.Range("AF65").Value = ''

